$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 173.9667156666667
$ws.Cells.Item(2, 8).Value = 521.9001470000001
$ws.Cells.Item(2, 9).Value = 0.339126905182122
$ws.Cells.Item(2, 10).Value = 0.339126905182122
$ws.Cells.Item(2, 13).Value = 58.95713633333333
$ws.Cells.Item(2, 14).Value = 176.871409
$ws.Cells.Item(2, 15).Value = 0.4863146960083892
$ws.Cells.Item(2, 16).Value = 0.4863146960083893
$ws.Cells.Item(2, 17).Value = 10256.5793730219
$ws.Cells.Item(2, 18).Value = 92309.21435719714
$ws.Cells.Item(2, 19).Value = 0.1649223978019095
$ws.Cells.Item(2, 20).Value = 0.1649223978019095
$ws.Cells.Item(3, 7).Value = 173.9667156666667
$ws.Cells.Item(3, 8).Value = 521.9001470000001
$ws.Cells.Item(3, 9).Value = 0.339126905182122
$ws.Cells.Item(3, 10).Value = 0.339126905182122
$ws.Cells.Item(3, 15).Value = 0.07416766570679004
$ws.Cells.Item(3, 16).Value = 0.07416766570679005
$ws.Cells.Item(3, 17).Value = 1564.226942918303
$ws.Cells.Item(3, 18).Value = 14078.04248626473
$ws.Cells.Item(3, 19).Value = 0.0251522509357259
$ws.Cells.Item(3, 20).Value = 0.02515225093572591
$ws.Cells.Item(4, 7).Value = 173.9667156666667
$ws.Cells.Item(4, 8).Value = 521.9001470000001
$ws.Cells.Item(4, 9).Value = 0.339126905182122
$ws.Cells.Item(4, 10).Value = 0.339126905182122
$ws.Cells.Item(4, 13).Value = 42.51661933333333
$ws.Cells.Item(4, 14).Value = 127.549858
$ws.Cells.Item(4, 15).Value = 0.3507032073181665
$ws.Cells.Item(4, 16).Value = 0.3507032073181665
$ws.Cells.Item(4, 17).Value = 7396.476626669903
$ws.Cells.Item(4, 18).Value = 66568.28964002914
$ws.Cells.Item(4, 19).Value = 0.1189328933352539
$ws.Cells.Item(4, 20).Value = 0.1189328933352539
$ws.Cells.Item(5, 7).Value = 173.9667156666667
$ws.Cells.Item(5, 8).Value = 521.9001470000001
$ws.Cells.Item(5, 9).Value = 0.339126905182122
$ws.Cells.Item(5, 10).Value = 0.339126905182122
$ws.Cells.Item(5, 13).Value = 10.76719366666667
$ws.Cells.Item(5, 14).Value = 32.301581
$ws.Cells.Item(5, 15).Value = 0.0888144309666542
$ws.Cells.Item(5, 16).Value = 0.08881443096665421
$ws.Cells.Item(5, 17).Value = 1873.133319136934
$ws.Cells.Item(5, 18).Value = 16858.19987223241
$ws.Cells.Item(5, 19).Value = 0.03011936310923266
$ws.Cells.Item(5, 20).Value = 0.03011936310923266
$ws.Cells.Item(6, 7).Value = 96.77942399999999
$ws.Cells.Item(6, 9).Value = 0.1886596894161923
$ws.Cells.Item(6, 10).Value = 0.1886596894161923
$ws.Cells.Item(6, 13).Value = 58.95713633333333
$ws.Cells.Item(6, 14).Value = 176.871409
$ws.Cells.Item(6, 15).Value = 0.4863146960083892
$ws.Cells.Item(6, 16).Value = 0.4863146960083893
$ws.Cells.Item(6, 17).Value = 5705.837695029471
$ws.Cells.Item(6, 18).Value = 51352.53925526524
$ws.Cells.Item(6, 19).Value = 0.0917479795074727
$ws.Cells.Item(6, 20).Value = 0.09174797950747271
$ws.Cells.Item(7, 7).Value = 96.77942399999999
$ws.Cells.Item(7, 9).Value = 0.1886596894161923
$ws.Cells.Item(7, 10).Value = 0.1886596894161923
$ws.Cells.Item(7, 15).Value = 0.07416766570679004
$ws.Cells.Item(7, 16).Value = 0.07416766570679005
$ws.Cells.Item(7, 17).Value = 870.1950942787198
$ws.Cells.Item(7, 18).Value = 7831.755848508478
$ws.Cells.Item(7, 19).Value = 0.01399244877696699
$ws.Cells.Item(7, 20).Value = 0.01399244877696699
$ws.Cells.Item(8, 7).Value = 96.77942399999999
$ws.Cells.Item(8, 9).Value = 0.1886596894161923
$ws.Cells.Item(8, 10).Value = 0.1886596894161923
$ws.Cells.Item(8, 13).Value = 42.51661933333333
$ws.Cells.Item(8, 14).Value = 127.549858
$ws.Cells.Item(8, 15).Value = 0.3507032073181665
$ws.Cells.Item(8, 16).Value = 0.3507032073181665
$ws.Cells.Item(8, 17).Value = 4114.733929507263
$ws.Cells.Item(8, 18).Value = 37032.60536556537
$ws.Cells.Item(8, 19).Value = 0.06616355816990779
$ws.Cells.Item(8, 20).Value = 0.06616355816990779
$ws.Cells.Item(9, 7).Value = 96.77942399999999
$ws.Cells.Item(9, 9).Value = 0.1886596894161923
$ws.Cells.Item(9, 10).Value = 0.1886596894161923
$ws.Cells.Item(9, 13).Value = 10.76719366666667
$ws.Cells.Item(9, 14).Value = 32.301581
$ws.Cells.Item(9, 15).Value = 0.0888144309666542
$ws.Cells.Item(9, 16).Value = 0.08881443096665421
$ws.Cells.Item(9, 17).Value = 1042.042801156448
$ws.Cells.Item(9, 18).Value = 9378.385210408031
$ws.Cells.Item(9, 19).Value = 0.01675570296184483
$ws.Cells.Item(9, 20).Value = 0.01675570296184484
$ws.Cells.Item(10, 7).Value = 120.8019793333333
$ws.Cells.Item(10, 8).Value = 362.405938
$ws.Cells.Item(10, 9).Value = 0.235488732624488
$ws.Cells.Item(10, 10).Value = 0.2354887326244879
$ws.Cells.Item(10, 13).Value = 58.95713633333333
$ws.Cells.Item(10, 14).Value = 176.871409
$ws.Cells.Item(10, 15).Value = 0.4863146960083892
$ws.Cells.Item(10, 16).Value = 0.4863146960083893
$ws.Cells.Item(10, 17).Value = 7122.138764891849
$ws.Cells.Item(10, 18).Value = 64099.24888402664
$ws.Cells.Item(10, 19).Value = 0.1145216314196787
$ws.Cells.Item(10, 20).Value = 0.1145216314196787
$ws.Cells.Item(11, 7).Value = 120.8019793333333
$ws.Cells.Item(11, 8).Value = 362.405938
$ws.Cells.Item(11, 9).Value = 0.235488732624488
$ws.Cells.Item(11, 10).Value = 0.2354887326244879
$ws.Cells.Item(11, 15).Value = 0.07416766570679004
$ws.Cells.Item(11, 16).Value = 0.07416766570679005
$ws.Cells.Item(11, 17).Value = 1086.194621235046
$ws.Cells.Item(11, 18).Value = 9775.751591115419
$ws.Cells.Item(11, 19).Value = 0.01746564959900868
$ws.Cells.Item(11, 20).Value = 0.01746564959900869
$ws.Cells.Item(12, 7).Value = 120.8019793333333
$ws.Cells.Item(12, 8).Value = 362.405938
$ws.Cells.Item(12, 9).Value = 0.235488732624488
$ws.Cells.Item(12, 10).Value = 0.2354887326244879
$ws.Cells.Item(12, 13).Value = 42.51661933333333
$ws.Cells.Item(12, 14).Value = 127.549858
$ws.Cells.Item(12, 15).Value = 0.3507032073181665
$ws.Cells.Item(12, 16).Value = 0.3507032073181665
$ws.Cells.Item(12, 17).Value = 5136.091770028534
$ws.Cells.Item(12, 18).Value = 46224.8259302568
$ws.Cells.Item(12, 19).Value = 0.08258665381869808
$ws.Cells.Item(12, 20).Value = 0.08258665381869808
$ws.Cells.Item(13, 7).Value = 120.8019793333333
$ws.Cells.Item(13, 8).Value = 362.405938
$ws.Cells.Item(13, 9).Value = 0.235488732624488
$ws.Cells.Item(13, 10).Value = 0.2354887326244879
$ws.Cells.Item(13, 13).Value = 10.76719366666667
$ws.Cells.Item(13, 14).Value = 32.301581
$ws.Cells.Item(13, 15).Value = 0.0888144309666542
$ws.Cells.Item(13, 16).Value = 0.08881443096665421
$ws.Cells.Item(13, 17).Value = 1300.698306798664
$ws.Cells.Item(13, 18).Value = 11706.28476118798
$ws.Cells.Item(13, 19).Value = 0.02091479778710248
$ws.Cells.Item(13, 20).Value = 0.02091479778710248
$ws.Cells.Item(14, 7).Value = 121.435997
$ws.Cells.Item(14, 8).Value = 364.307991
$ws.Cells.Item(14, 9).Value = 0.2367246727771976
$ws.Cells.Item(14, 10).Value = 0.2367246727771976
$ws.Cells.Item(14, 13).Value = 58.95713633333333
$ws.Cells.Item(14, 14).Value = 176.871409
$ws.Cells.Item(14, 15).Value = 0.4863146960083892
$ws.Cells.Item(14, 16).Value = 0.4863146960083893
$ws.Cells.Item(14, 17).Value = 7159.518630903257
$ws.Cells.Item(14, 18).Value = 64435.66767812932
$ws.Cells.Item(14, 19).Value = 0.1151226872793283
$ws.Cells.Item(14, 20).Value = 0.1151226872793283
$ws.Cells.Item(15, 7).Value = 121.435997
$ws.Cells.Item(15, 8).Value = 364.307991
$ws.Cells.Item(15, 9).Value = 0.2367246727771976
$ws.Cells.Item(15, 10).Value = 0.2367246727771976
$ws.Cells.Item(15, 15).Value = 0.07416766570679004
$ws.Cells.Item(15, 16).Value = 0.07416766570679005
$ws.Cells.Item(15, 17).Value = 1091.89541010541
$ws.Cells.Item(15, 18).Value = 9827.05869094869
$ws.Cells.Item(15, 19).Value = 0.01755731639508845
$ws.Cells.Item(15, 20).Value = 0.01755731639508845
$ws.Cells.Item(16, 7).Value = 121.435997
$ws.Cells.Item(16, 8).Value = 364.307991
$ws.Cells.Item(16, 9).Value = 0.2367246727771976
$ws.Cells.Item(16, 10).Value = 0.2367246727771976
$ws.Cells.Item(16, 13).Value = 42.51661933333333
$ws.Cells.Item(16, 14).Value = 127.549858
$ws.Cells.Item(16, 15).Value = 0.3507032073181665
$ws.Cells.Item(16, 16).Value = 0.3507032073181665
$ws.Cells.Item(16, 17).Value = 5163.048057812808
$ws.Cells.Item(16, 18).Value = 46467.43252031528
$ws.Cells.Item(16, 19).Value = 0.08302010199430665
$ws.Cells.Item(16, 20).Value = 0.08302010199430666
$ws.Cells.Item(17, 7).Value = 121.435997
$ws.Cells.Item(17, 8).Value = 364.307991
$ws.Cells.Item(17, 9).Value = 0.2367246727771976
$ws.Cells.Item(17, 10).Value = 0.2367246727771976
$ws.Cells.Item(17, 13).Value = 10.76719366666667
$ws.Cells.Item(17, 14).Value = 32.301581
$ws.Cells.Item(17, 15).Value = 0.0888144309666542
$ws.Cells.Item(17, 16).Value = 0.08881443096665421
$ws.Cells.Item(17, 17).Value = 1307.524897803752
$ws.Cells.Item(17, 18).Value = 11767.72408023377
$ws.Cells.Item(17, 19).Value = 0.02102456710847422
$ws.Cells.Item(17, 20).Value = 0.02102456710847422
